# Remove the two obsolete "LatestRun_SkapaReskort_*" asset rows from the
# Assets sheet (rows 3 and 4). Deleting the entire rows shifts everything
# below up by two rows and drops the two now-unused shared strings from
# the workbook's shared string table (which also renumbers every other
# <c t="s"> reference, including the one on the Settings sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

$ws.Rows("3:4").Delete()

# Mirror the author's final selection: the whole of (new) row 3 - i.e.
# what used to be row 5 (FilePath_Masterfile_Elevresor) - ends up selected.
$ws.Range("3:3").Select() | Out-Null
